# Add two new header columns (I: "I0", J: "IF") with the same formatting
# as the existing header cells, plus their data values for rows 2 and 3.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting of the last existing header cell (H1) onto the new
# header cells before writing their text, so the new cells pick up the
# bold/centered/bordered header style (same style index as the other
# header cells).
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

$ws.Range("I2").Value = 11
$ws.Range("J2").Value = 11

$ws.Range("I3").Value = 8
$ws.Range("J3").Value = 8
